$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the emoji status icons in column A (statut) with the new symbols.
# 📘 -> ⚠️  and  📗 -> ✅
$ws.Range("A2").Value = "⚠️"
$ws.Range("A3").Value = "✅"
$ws.Range("A4").Value = "✅"
